# "this is the task of section OOS"
#
# Adds a second row of submission data (name / email / repo link) to the
# existing "name | email | Repo Link" header row on Sheet1, turns the email
# cell into a mailto: hyperlink (Excel's built-in "Hyperlink" cell style),
# and leaves the selection on C12 (matching the author's last on-screen
# selection when the file was saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row under the header (row 1): name, email, repo link.
$ws.Range("A2").Value = "كيرلس جمال نجيب"
$ws.Range("B2").Value = "kerolosgnaguib@gmail.com"
$ws.Range("C2").Value = "https://github.com/harounwaka125/OOS-Project"

# Turn the email cell into a live mailto: hyperlink. Excel auto-applies the
# built-in "Hyperlink" style (underline + theme color) to the cell; since
# B2 already holds the display text, we don't pass TextToDisplay so the
# existing cell text/string is kept as-is.
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:kerolosgnaguib@gmail.com") | Out-Null

# Match the saved selection/active cell.
$ws.Range("C12").Select() | Out-Null
